$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $ok = $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find failed for: $find"
    }
}

# --- Paragraph: "The country you do live in affects ..." (cause-of-death reflection) ---
Replace-Text "the likely chance of you dying from certain death" `
             "the likely chance of you dying from a certain death"
Replace-Text "one we see everyday is Covid-19, it is very big in the United States" `
             "Covid-19, it is very big in the United States"
Replace-Text "from covid. That is how probably is in general though. " `
             "from covid. There is causes such as Malaria that are only common in certain parts of the world no matter what. "

# --- Paragraph: " Outcome of your EDA" header (merge runs, text unchanged) ---
Replace-Text " Outcome of your EDA" " Outcome of your EDA"

# --- Paragraph: "I believe the outcome of my EDA ..." ---
Replace-Text "I believe the outcome of my EDA was just insightfulness on what I learned the whole semester. " `
             "I believe the outcome of my EDA was very beneficial and insightful. "

# --- Paragraph: "I feel like I did miss a lot in my analysis, ..." ---
Replace-Text "I feel like I did miss a lot in my analysis, but it was the different languages that made it harder for me because in the book they had their own code versus when I started using my own data, I had to find different dictionaries that support python. That’s not all bad though because I learned a lot more during this project as well as it reinforced a lot. " `
             "I feel like I missed some libraries that would have helped me describe my data better but I feel like I will learn that soon in the next upcoming classes."

# --- Paragraph: "I think the professor did a good job ..." ---
Replace-Text "did a good job at listing" "did a good job at listing"
Replace-Text "them." `
             "them, there were some from newer libraries that I came across that we hadn’t learned yet so I plan to look into that shortly after I finish this semester. "

# --- Paragraph: "I feel like my data could have been better ..." ---
Replace-Text "could have been labeled because" "could have been labeled because"
Replace-Text "something is off as if it could have been better. I also think I just did not choose a topic with a lot of variables which eventually bit me in the butt, but it was too late to change it, so I made it work. " `
             "something is off because there could have been a different more randomized sorting that doesn’t relate to the alphabetical name. I also think I just did not choose a topic with enough different dependent variables. "

# --- Paragraph: "One main challenge I faced ..." ---
Replace-Text " It is some other stuff that haven’t stuck yet but I hope to eventually get there. " ""
